# Scoreboard.xlsx edit: add two "semi-final" sheets (SFM / SFF) holding
# Snatch / Clean-and-Jerk lift numbers for the qualifying men's and
# women's teams, positioned right before the ScoreMatrix sheet.

$wb = $excel.ActiveWorkbook

$scoreMatrix = $wb.Worksheets.Item("ScoreMatrix")

# Insert "SFM" right before ScoreMatrix (after ScoreF).
$sfm = $wb.Worksheets.Add($scoreMatrix)
$sfm.Name = "SFM"

# Insert "SFF" right after the freshly created SFM (still before ScoreMatrix).
$sff = $wb.Worksheets.Add($null, $sfm)
$sff.Name = "SFF"

# ---- SFM: men's semi-final lift data -------------------------------------
$sfm.Cells.Item(1,1).Value = "Team"
$sfm.Cells.Item(1,2).Value = "Snatch"
$sfm.Cells.Item(1,3).Value = "Clean and Jerk"

$sfmRows = @(
    @("Håkon Konningen og Njål Christensen", 75, 100),
    @("Anders Magnus Nes og Anders Vinnes Jacobsen", 78, 111),
    @("Andreas Melheim Hansen og Jørgen Skarsmo", 80, 99),
    @("Eirik Berge og Julian Haug", 66, 85),
    @("Kasper Støen Nerbøvik og Håvard Idland", 52, 98),
    @("Magnus Øslebye og Vegard Tangen", 77, 104),
    @("Anders J. Svalestuen og Gabriel Kristiansen", 65, 101)
)

$r = 2
foreach ($row in $sfmRows) {
    $sfm.Cells.Item($r,1).Value = $row[0]
    $sfm.Cells.Item($r,2).Value = $row[1]
    $sfm.Cells.Item($r,3).Value = $row[2]
    $r = $r + 1
}

$sfm.Columns.Item(1).AutoFit() | Out-Null
$sfm.Columns.Item(3).AutoFit() | Out-Null
$sfm.Range("C9").Select() | Out-Null

# ---- SFF: women's semi-final lift data ------------------------------------
$sff.Cells.Item(1,1).Value = "Team"
$sff.Cells.Item(1,2).Value = "Snatch"
$sff.Cells.Item(1,3).Value = "Clean and Jerk"

$sffRows = @(
    @("Renate Berntsen Hansen og Karoline Granås", 66, 90),
    @("Maria Hanssen og Cecilie Rabben", 59, 85),
    @("Victoria Christensen og Helene Rye Martinsen", 72, 81),
    @("Marianne U. Henriksen og Mari S. Andersen", 78, 75),
    @("Dawn Stewart og Marie Vik", 55, 91),
    @("Sara Yuzer og Martine Baalsrud", 69, 80),
    @("Frid Kaspersen og Renate Loraas", 70, 78)
)

$r = 2
foreach ($row in $sffRows) {
    $sff.Cells.Item($r,1).Value = $row[0]
    $sff.Cells.Item($r,2).Value = $row[1]
    $sff.Cells.Item($r,3).Value = $row[2]
    $r = $r + 1
}

$sff.Columns.Item(1).AutoFit() | Out-Null
$sff.Columns.Item(3).AutoFit() | Out-Null
$sff.Range("J10").Select() | Out-Null

# Restore the originally active tab (ScoreF) so the workbook still opens
# on the same sheet it did before the two inserts.
$wb.Worksheets.Item("ScoreF").Activate()
